$d = $word.ActiveDocument
$sub = $d.Range(18, 23)
Write-Output "sub text: [$($sub.Text)]"
$sub.Font.Bold = 1
Write-Output "done"
